$wb = $excel.ActiveWorkbook

# --- "About" sheet: update the "date updated" field (C1) ---
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("C1").Value = 45392

# --- "MCF" sheet: update capacity factor values from 0.85/0.95 to 1 ---
$wsMCF = $wb.Worksheets.Item("MCF")

$cellsToUpdate = @("B2","B3","B4","B6","B10","B11","B12","B13","B14","B16","B17","B18")
foreach ($cell in $cellsToUpdate) {
    $wsMCF.Range($cell).Value = 1
}

# Move the active selection on the MCF sheet to B17 (matches final saved selection)
$wsMCF.Activate()
$wsMCF.Range("B17").Select()
